$d = $word.ActiveDocument

# Fix contact information missing from short resumes: insert a centered
# contact-info line directly below the "Dheeraj Chand" title, matching the
# long-resume layout. Using Find/Replace with a "^p" paragraph-mark token in
# the replacement text splits "Dheeraj Chand" into its own paragraph and adds
# a new paragraph after it (inheriting the same centered alignment) without
# also inheriting the title run's bold/large-size character formatting.
$d.Content.Find.Execute(
    "Dheeraj Chand", $true, $false, $false, $false, $false, $true, 1, $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
